$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.85
$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 2.02
$ws.Range("I2").Value = 2.22
$ws.Range("J2").Value = 3.3
$ws.Range("K2").Value = 3.75
$ws.Range("L2").Value = 1.45
$ws.Range("O2").Value = 1.38
$ws.Range("P2").Value = 1.74
$ws.Range("R2").Value = 1.28
$ws.Range("V2").Value = 1.82
$ws.Range("W2").Value = 1.29
$ws.Range("Z2").Value = 13.5
$ws.Range("AA2").Value = 34
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 32
$ws.Range("AK2").Value = 60
$ws.Range("AO2").Value = 21

# Row 3
$ws.Range("G3").Value = 1.98
$ws.Range("H3").Value = 4.3
$ws.Range("K3").Value = 3.95
$ws.Range("N3").Value = 3.6
$ws.Range("Q3").Value = 1.93
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 38
$ws.Range("AB3").Value = 9.4
$ws.Range("AC3").Value = 8.4
$ws.Range("AG3").Value = 1000
$ws.Range("AN3").Value = 17

# Row 4
$ws.Range("N4").Value = 3.05
$ws.Range("O4").Value = 1.38
$ws.Range("P4").Value = 1.71
$ws.Range("Q4").Value = 2.14
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.98

# Row 5
$ws.Range("F5").Value = 1.48
$ws.Range("G5").Value = 1.56
$ws.Range("H5").Value = 6.6
$ws.Range("I5").Value = 8.6
$ws.Range("J5").Value = 3.95
$ws.Range("K5").Value = 5.1
$ws.Range("N5").Value = 3.3
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 1.35
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 1.84
$ws.Range("AF5").Value = 1000

# Row 6
$ws.Range("G6").Value = 3.4
$ws.Range("H6").Value = 2.34
$ws.Range("I6").Value = 2.54
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 3.65
$ws.Range("L6").Value = 1.35
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 1.96
$ws.Range("Q6").Value = 1.92
$ws.Range("U6").Value = 2.28
$ws.Range("V6").Value = 1.64
$ws.Range("W6").Value = 1.42
$ws.Range("X6").Value = 18

# Row 7
$ws.Range("G7").Value = 1.44
$ws.Range("J7").Value = 5
$ws.Range("L7").Value = 1.25
$ws.Range("N7").Value = 4.9
$ws.Range("P7").Value = 2.32
$ws.Range("T7").Value = 1.88
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.1
$ws.Range("W7").Value = 3.25
$ws.Range("Z7").Value = 110
$ws.Range("AK7").Value = 17.5

# Row 8
$ws.Range("F8").Value = 1.66
$ws.Range("N8").Value = 2.48
$ws.Range("P8").Value = 1.5
$ws.Range("R8").Value = 1.16
$ws.Range("U8").Value = 1.47
$ws.Range("X8").Value = 990
$ws.Range("Y8").Value = 990
$ws.Range("AC8").Value = 990
$ws.Range("AF8").Value = 980
$ws.Range("AG8").Value = 990
